$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.934.87'
$ws.Range("E2").Value = '  -0.30%  '

$ws.Range("D3").Value = '1.910.34'
$ws.Range("E3").Value = '  -0.05%  '

$ws.Range("D4").Value = '0.9983'
$ws.Range("E4").Value = '  -0.74%  '

$ws.Range("D5").Value = '313.03'

$ws.Range("D6").Value = '0.9983'
$ws.Range("E6").Value = '  -0.64%  '

$ws.Range("D7").Value = '0.5013'

$ws.Range("D8").Value = '0.3816'
$ws.Range("E8").Value = '  +0.10%  '

$ws.Range("D9").Value = '0.07307'
$ws.Range("E9").Value = '  -0.74%  '

$ws.Range("E10").Value = '  -2.40%  '

$ws.Range("D11").Value = '21.26'

$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.941.22'
$ws.Range("E12").Value = '  +1.67%  '

$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").Value = '0.07669'
$ws.Range("E13").Value = '  -1.86%  '

$ws.Range("D14").Value = '5.485'
$ws.Range("E14").Value = '  -0.48%  '

$ws.Range("D15").Value = '92.82'
$ws.Range("E15").Value = '  +0.80%  '

$ws.Range("D16").Value = '0.9988'
$ws.Range("E16").Value = '  -0.69%  '

$ws.Range("D17").Value = '0.000008743'
$ws.Range("E17").Value = '  -1.45%  '

$ws.Range("D18").Value = '0.9981'
$ws.Range("E18").Value = '  -0.65%  '

$ws.Range("D19").Value = '27.971.10'
$ws.Range("E19").Value = '  -0.28%  '

$ws.Range("D20").Value = '14.68'
$ws.Range("E20").Value = '  -0.44%  '

$ws.Range("D21").Value = '5.188'
$ws.Range("E21").Value = '  +0.28%  '

$ws.Range("D22").Value = '2.154.11'
$ws.Range("E22").Value = '  +0.85%  '

$ws.Range("E23").Value = '  -0.45%  '

$ws.Range("D24").Value = '6.615'
$ws.Range("E24").Value = '  -0.34%  '

$ws.Range("D25").Value = '153.13'
$ws.Range("E25").Value = '  -2.58%  '

$ws.Range("D26").Value = '1.842'
$ws.Range("E26").Value = '  -4.00%  '

$ws.Range("E27").Value = '  +3.35%  '

$ws.Range("D28").Value = '18.43'

$ws.Range("D29").Value = '115.50'
$ws.Range("E29").Value = '  -1.36%  '

$ws.Range("D30").Value = '4.934'
$ws.Range("E30").Value = '  -0.79%  '

$ws.Range("D31").Value = '0.09026'
$ws.Range("E31").Value = '  +0.82%  '

$ws.Range("D32").Value = '3.210'
$ws.Range("E32").Value = '  -2.36%  '

$ws.Range("D33").Value = '4.852'
$ws.Range("E33").Value = '  +3.99%  '

$ws.Range("D34").Value = '1.239'
$ws.Range("E34").Value = '  -1.43%  '

$ws.Range("D35").Value = '0.7818'
$ws.Range("E35").Value = '  +0.92%  '

$ws.Range("B36").Value = 'VeChain'
$ws.Range("C36").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D36").Value = '0.02088'
$ws.Range("E36").Value = '  +1.92%  '

$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D37").Value = '2.603'
$ws.Range("E37").Value = '  -0.33%  '

$ws.Range("D38").Value = '3.069'
$ws.Range("E38").Value = '  +2.39%  '

$ws.Range("E39").Value = '  -1.33%  '

$ws.Range("D40").Value = '0.5558'
$ws.Range("E40").Value = '  +0.65%  '

$ws.Range("D41").Value = '0.05286'

$ws.Range("D42").Value = '6.890'
$ws.Range("E42").Value = '  -2.16%  '

$ws.Range("D43").Value = '113.71'
$ws.Range("E43").Value = '  +4.67%  '

$ws.Range("D44").Value = '8.550'
$ws.Range("E44").Value = '  +0.42%  '

$ws.Range("D45").Value = '0.1519'
$ws.Range("E45").Value = '  -0.66%  '

$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = '10.64'
$ws.Range("E46").Value = '  -0.79%  '

$ws.Range("B47").Value = 'Decentraland'
$ws.Range("C47").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D47").Value = '0.4836'
$ws.Range("E47").Value = '  +0.06%  '

$ws.Range("D48").Value = '0.9974'
$ws.Range("E48").Value = '  -0.78%  '

$ws.Range("D49").Value = '1.642'
$ws.Range("E49").Value = '  -0.51%  '

$ws.Range("D50").Value = '67.63'
$ws.Range("E50").Value = '  -0.66%  '

$ws.Range("E51").Value = '  -0.52%  '
